$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 4
$ws.Range("H4").Value = 284.25
$ws.Range("I4").Value = 320.57144
$ws.Range("J4").Value = 30
$ws.Range("K4").Value = 320.57144
$ws.Range("L4").Value = 30
$ws.Range("M4").Value = -206.57144
$ws.Range("N4").Value = -258

# Row 81
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996

# Row 84
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984

# Row 86
$ws.Range("H86").Value = 76965720
$ws.Range("I86").Value = 3420.5
$ws.Range("K86").Value = 3420.5
$ws.Range("M86").Value = -2297.5

# Row 89
$ws.Range("H89").Value = 76965720
$ws.Range("I89").Value = 3420.5
$ws.Range("K89").Value = 17102.5
$ws.Range("M89").Value = -11486.5

# Row 98
$ws.Range("H98").Value = 25709.912
$ws.Range("I98").Value = 26514.863
$ws.Range("J98").Value = 8001
$ws.Range("K98").Value = 26514.863
$ws.Range("L98").Value = 8001
$ws.Range("M98").Value = -25016.863
$ws.Range("N98").Value = -10997

# Row 106
$ws.Range("H106").Value = 5617223.5
$ws.Range("I106").Value = 6863728.5
$ws.Range("K106").Value = 6863728.5
$ws.Range("M106").Value = -6863097.5

# Row 112
$ws.Range("H112").Value = 2665.6667
$ws.Range("J112").Value = 2043.1111
$ws.Range("L112").Value = 6129.3333
$ws.Range("N112").Value = -8345.3333

# Row 122
$ws.Range("H122").Value = 25709.912
$ws.Range("I122").Value = 26514.863
$ws.Range("J122").Value = 8001
$ws.Range("K122").Value = 79544.58900000001
$ws.Range("L122").Value = 24003
$ws.Range("M122").Value = -77094.58900000001
$ws.Range("N122").Value = -28903

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 8533.064
$ws.Range("I32").Value = 8533.064
$ws.Range("K32").Value = 8533.064
$ws.Range("M32").Value = -8246.064

# Row 61
$ws.Range("H61").Value = 6826.7
$ws.Range("J61").Value = 5230.0835
$ws.Range("L61").Value = 5230.0835
$ws.Range("N61").Value = -5654.0835

# Row 82
$ws.Range("H82").Value = 64724.5
$ws.Range("J82").Value = 64724.5
$ws.Range("L82").Value = 64724.5
$ws.Range("N82").Value = -65446.5

# Row 85
$ws.Range("H85").Value = 64724.5
$ws.Range("J85").Value = 64724.5
$ws.Range("L85").Value = 64724.5
$ws.Range("N85").Value = -67220.5

# Row 132
$ws.Range("H132").Value = 5188.5
$ws.Range("I132").Value = 5066.6665
$ws.Range("K132").Value = 15199.9995
$ws.Range("M132").Value = -12669.9995

# Row 136
$ws.Range("H136").Value = 6826.7
$ws.Range("J136").Value = 5230.0835
$ws.Range("L136").Value = 15690.2505
$ws.Range("N136").Value = -20790.2505

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 11
$ws.Range("H11").Value = 500
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 0
$ws.Range("N11").Value = -780
$ws.Range("M11").ClearContents()

# Row 16
$ws.Range("H16").Value = 2709.6428
$ws.Range("I16").Value = 2798.125
$ws.Range("K16").Value = 2798.125
$ws.Range("M16").Value = -2511.125

# Row 113
$ws.Range("H113").Value = 2709.6428
$ws.Range("I113").Value = 2798.125
$ws.Range("K113").Value = 2798.125
$ws.Range("M113").Value = -628.125

# Row 134
$ws.Range("H134").Value = 1564.5667
$ws.Range("I134").Value = 1216.8422
$ws.Range("K134").Value = 3650.5266
$ws.Range("M134").Value = -1115.5266

# Row 141
$ws.Range("H141").Value = 312965.28
$ws.Range("J141").Value = 360148.28
$ws.Range("L141").Value = 360148.28
$ws.Range("N141").Value = -370508.28

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 23
$ws.Range("H23").Value = 13889779
$ws.Range("I23").Value = 883.3333
$ws.Range("J23").Value = 18519410
$ws.Range("K23").Value = 2649.9999
$ws.Range("L23").Value = 55558230
$ws.Range("M23").Value = -2414.9999
$ws.Range("N23").Value = -55558700

# Row 137
$ws.Range("H137").Value = 10800.65
$ws.Range("I137").Value = 3799.8
$ws.Range("J137").Value = 13134.267
$ws.Range("K137").Value = 11399.4
$ws.Range("L137").Value = 39402.801
$ws.Range("M137").Value = -6299.400000000001
$ws.Range("N137").Value = -49602.801

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 4
$ws.Range("H4").Value = 6450
$ws.Range("I4").Value = 6933.3335
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 6933.3335
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -6821.3335
$ws.Range("N4").Value = -5224

# Row 97
$ws.Range("H97").Value = 8746.786
$ws.Range("I97").Value = 10850.546
$ws.Range("J97").Value = 1033
$ws.Range("K97").Value = 10850.546
$ws.Range("L97").Value = 1033
$ws.Range("M97").Value = -10354.546
$ws.Range("N97").Value = -2025

# Row 107
$ws.Range("H107").Value = 774.94116
$ws.Range("J107").Value = 785
$ws.Range("L107").Value = 785
$ws.Range("N107").Value = -4625

# Row 132
$ws.Range("H132").Value = 3469.3823
$ws.Range("I132").Value = 3340.1667
$ws.Range("K132").Value = 10020.5001
$ws.Range("M132").Value = -7490.500100000001

# Row 141
$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 16
$ws.Range("H16").Value = 8596.842000000001
$ws.Range("I16").Value = 9083.625
$ws.Range("J16").Value = 6000.6665
$ws.Range("K16").Value = 9083.625
$ws.Range("L16").Value = 6000.6665
$ws.Range("M16").Value = -8913.625
$ws.Range("N16").Value = -6340.6665

# Row 46
$ws.Range("H46").Value = 3666.5557
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612

# Row 68
$ws.Range("H68").Value = 5729.3687
$ws.Range("I68").Value = 3042.7144
$ws.Range("J68").Value = 7296.5835
$ws.Range("K68").Value = 3042.7144
$ws.Range("L68").Value = 7296.5835
$ws.Range("M68").Value = -2293.7144
$ws.Range("N68").Value = -8794.583500000001

# Row 71
$ws.Range("H71").Value = 5729.3687
$ws.Range("I71").Value = 3042.7144
$ws.Range("J71").Value = 7296.5835
$ws.Range("K71").Value = 15213.572
$ws.Range("L71").Value = 36482.9175
$ws.Range("M71").Value = -11469.572
$ws.Range("N71").Value = -43970.9175

# Row 100
$ws.Range("H100").Value = 7154.273
$ws.Range("I100").Value = 5949.6665
$ws.Range("K100").Value = 5949.6665
$ws.Range("M100").Value = -5408.6665

# Row 101
$ws.Range("H101").Value = 26667
$ws.Range("J101").Value = 26667
$ws.Range("L101").Value = 26667
$ws.Range("N101").Value = -33157

# Row 122
$ws.Range("H122").Value = 7491.154
$ws.Range("I122").Value = 8042.778
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 24128.334
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -21678.334
$ws.Range("N122").Value = -23650

# Row 132
$ws.Range("H132").Value = 682575.25
$ws.Range("I132").Value = 998527.2
$ws.Range("J132").Value = 5535.4287
$ws.Range("K132").Value = 2995581.6
$ws.Range("L132").Value = 16606.2861
$ws.Range("M132").Value = -2993051.6
$ws.Range("N132").Value = -21666.2861

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 68
$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26622

# Row 71
$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -83112

# Row 98
$ws.Range("H98").Value = 53999.668
$ws.Range("J98").Value = 53999.668
$ws.Range("L98").Value = 53999.668
$ws.Range("N98").Value = -59989.668

# Row 100
$ws.Range("H100").Value = 29157.846
$ws.Range("I100").Value = 6339.4443
$ws.Range("J100").Value = 80499.25
$ws.Range("K100").Value = 12678.8886
$ws.Range("L100").Value = 160998.5
$ws.Range("M100").Value = -12137.8886
$ws.Range("N100").Value = -162080.5

# Row 132
$ws.Range("H132").Value = 4922.0786
$ws.Range("I132").Value = 5109.528
$ws.Range("J132").Value = 4472.2
$ws.Range("K132").Value = 15328.584
$ws.Range("L132").Value = 13416.6
$ws.Range("M132").Value = -12798.584
$ws.Range("N132").Value = -18476.6
